$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 5843.6924
$ws.Range("I21").Value = 4537.3335
$ws.Range("J21").Value = 6963.4287
$ws.Range("K21").Value = 4537.3335
$ws.Range("L21").Value = 6963.4287
$ws.Range("M21").Value = -4069.3335
$ws.Range("N21").Value = -7899.4287
$ws.Range("H23").Value = 5843.6924
$ws.Range("I23").Value = 4537.3335
$ws.Range("J23").Value = 6963.4287
$ws.Range("K23").Value = 4537.3335
$ws.Range("L23").Value = 6963.4287
$ws.Range("M23").Value = -4303.3335
$ws.Range("N23").Value = -7431.4287
$ws.Range("H28").Value = 1190.6364
$ws.Range("J28").Value = 1800
$ws.Range("L28").Value = 1800
$ws.Range("N28").Value = -2770
$ws.Range("H29").Value = 753.3
$ws.Range("I29").Value = 505.5
$ws.Range("J29").Value = 1125
$ws.Range("K29").Value = 1516.5
$ws.Range("L29").Value = 3375
$ws.Range("M29").Value = -1235.5
$ws.Range("N29").Value = -3937
$ws.Range("H38").Value = 1792594.4
$ws.Range("J38").Value = 755
$ws.Range("L38").Value = 2265
$ws.Range("N38").Value = -3009
$ws.Range("H58").Value = 489768.06
$ws.Range("I58").Value = 797674.2
$ws.Range("J58").Value = 2250
$ws.Range("K58").Value = 2393022.6
$ws.Range("L58").Value = 6750
$ws.Range("M58").Value = -2392872.6
$ws.Range("N58").Value = -7050
$ws.Range("H80").Value = 106278.21
$ws.Range("I80").Value = 200860.4
$ws.Range("J80").Value = 72498.86
$ws.Range("K80").Value = 602581.2
$ws.Range("L80").Value = 217496.58
$ws.Range("M80").Value = -601583.2
$ws.Range("N80").Value = -219492.58
$ws.Range("H83").Value = 106278.21
$ws.Range("I83").Value = 200860.4
$ws.Range("J83").Value = 72498.86
$ws.Range("K83").Value = 1807743.6
$ws.Range("L83").Value = 652489.74
$ws.Range("M83").Value = -1802751.6
$ws.Range("N83").Value = -662473.74
$ws.Range("H88").Value = 5687.5
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 6285.7144
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 6285.7144
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -7097.7144
$ws.Range("H91").Value = 5687.5
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 6285.7144
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 6285.7144
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -9093.714400000001
$ws.Range("H92").Value = 515.3889
$ws.Range("I92").Value = 515.3889
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 515.3889
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 732.6111
$ws.Range("N92").ClearContents()
$ws.Range("H101").Value = 462.85715
$ws.Range("I101").Value = 462.85715
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1388.57145
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 233.4285500000001
$ws.Range("N101").ClearContents()
$ws.Range("H112").Value = 1320
$ws.Range("I112").Value = 650
$ws.Range("J112").Value = 1766.6666
$ws.Range("K112").Value = 1950
$ws.Range("L112").Value = 5299.9998
$ws.Range("M112").Value = -842
$ws.Range("N112").Value = -7515.9998
$ws.Range("H113").Value = 57772.445
$ws.Range("J113").Value = 2272.6365
$ws.Range("L113").Value = 2272.6365
$ws.Range("N113").Value = -8780.636500000001
$ws.Range("H127").Value = 2403.465
$ws.Range("J127").Value = 2556.025
$ws.Range("L127").Value = 7668.075000000001
$ws.Range("N127").Value = -17588.075
$ws.Range("H129").Value = 807.1579
$ws.Range("J129").Value = 1127.7
$ws.Range("L129").Value = 3383.1
$ws.Range("N129").Value = -13383.1
$ws.Range("H132").Value = 5213257.5
$ws.Range("I132").Value = 5439616.5
$ws.Range("K132").Value = 16318849.5
$ws.Range("M132").Value = -16316319.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38550.207
$ws.Range("I32").Value = 7168.3423
$ws.Range("K32").Value = 7168.3423
$ws.Range("M32").Value = -6881.3423

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 243.4375
$ws.Range("I64").Value = 85
$ws.Range("J64").Value = 296.25
$ws.Range("K64").Value = 85
$ws.Range("L64").Value = 296.25
$ws.Range("M64").Value = 140
$ws.Range("N64").Value = -746.25
$ws.Range("H67").Value = 243.4375
$ws.Range("I67").Value = 85
$ws.Range("J67").Value = 296.25
$ws.Range("K67").Value = 85
$ws.Range("L67").Value = 296.25
$ws.Range("M67").Value = 695
$ws.Range("N67").Value = -1856.25
$ws.Range("H86").Value = 101874.45
$ws.Range("J86").Value = 3201.75
$ws.Range("L86").Value = 3201.75
$ws.Range("N86").Value = -5447.75
$ws.Range("H89").Value = 101874.45
$ws.Range("J89").Value = 3201.75
$ws.Range("L89").Value = 16008.75
$ws.Range("N89").Value = -27240.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 18832.334
$ws.Range("I99").Value = 2560
$ws.Range("J99").Value = 35104.668
$ws.Range("K99").Value = 2560
$ws.Range("L99").Value = 35104.668
$ws.Range("M99").Value = -1062
$ws.Range("N99").Value = -38100.668
$ws.Range("H126").Value = 18832.334
$ws.Range("I126").Value = 2560
$ws.Range("J126").Value = 35104.668
$ws.Range("K126").Value = 7680
$ws.Range("L126").Value = 105314.004
$ws.Range("M126").Value = -5210
$ws.Range("N126").Value = -110254.004
$ws.Range("H132").Value = 45458190
$ws.Range("I132").Value = 45458236
$ws.Range("K132").Value = 136374708
$ws.Range("M132").Value = -136372178

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 1666.6666
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 3000
$ws.Range("M125").Value = 1920
$ws.Range("H140").Value = 1419.3103
$ws.Range("I140").Value = 1017.1429
$ws.Range("K140").Value = 3051.4287
$ws.Range("M140").Value = 2128.5713

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1654.48
$ws.Range("I80").Value = 1452.2
$ws.Range("J80").Value = 1957.9
$ws.Range("K80").Value = 1452.2
$ws.Range("L80").Value = 1957.9
$ws.Range("M80").Value = -454.2
$ws.Range("N80").Value = -3953.9
$ws.Range("H83").Value = 1654.48
$ws.Range("I83").Value = 1452.2
$ws.Range("J83").Value = 1957.9
$ws.Range("K83").Value = 7261
$ws.Range("L83").Value = 9789.5
$ws.Range("M83").Value = -2269
$ws.Range("N83").Value = -19773.5
$ws.Range("H132").Value = 2894
$ws.Range("I132").Value = 2283.3333
$ws.Range("J132").Value = 3260.4
$ws.Range("K132").Value = 6849.999899999999
$ws.Range("L132").Value = 9781.200000000001
$ws.Range("M132").Value = -4319.999899999999
$ws.Range("N132").Value = -14841.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4547954.5
$ws.Range("I62").Value = 50000000
$ws.Range("J62").Value = 2750
$ws.Range("K62").Value = 50000000
$ws.Range("L62").Value = 2750
$ws.Range("M62").Value = -49999376
$ws.Range("N62").Value = -3998
$ws.Range("H65").Value = 4547954.5
$ws.Range("I65").Value = 50000000
$ws.Range("J65").Value = 2750
$ws.Range("K65").Value = 250000000
$ws.Range("L65").Value = 13750
$ws.Range("M65").Value = -249996880
$ws.Range("N65").Value = -19990
$ws.Range("H81").Value = 500984.75
$ws.Range("I81").Value = 500480
$ws.Range("J81").Value = 501489.5
$ws.Range("K81").Value = 1000960
$ws.Range("L81").Value = 1002979
$ws.Range("M81").Value = -999899
$ws.Range("N81").Value = -1005101
$ws.Range("H84").Value = 500984.75
$ws.Range("I84").Value = 500480
$ws.Range("J84").Value = 501489.5
$ws.Range("K84").Value = 5004800
$ws.Range("L84").Value = 5014895
$ws.Range("M84").Value = -4999496
$ws.Range("N84").Value = -5025503

Write-Host "Applied all profit-table updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets."
